# Rebuild the table: drop the "Sourabh" contributor entirely, reorder the
# remaining contributors so Ritesh comes right after Anshuman (before
# Rohan), fix the pistachios item-name typo, drop the "Tortilla" row, and
# recompute the contribution/value columns so the price of every item is
# split evenly across the four remaining contributors.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the old "Sourabh" columns (K, L) completely so no stale data /
# formatting remains once the table shrinks to 10 columns.
$ws.Columns.Item(11).Delete()
$ws.Columns.Item(11).Delete()

# Remove the last data row ("Tortilla").
$ws.Rows.Item(10).Delete()

$header = @("item_name", "price", "Anshuman_cont", "Ritesh_cont", "Rohan_cont", "Ashish_cont", "Anshuman_value", "Ritesh_value", "Rohan_value", "Ashish_value")
for ($c = 1; $c -le $header.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $header[$c - 1]
}

$items = @(
    @("45612 MILK WHOLE ", 1.55),
    @("22026 STIR FRY LARGE", 2.58),
    @("810873 FAMILY PACK TOMATO", 1.39),
    @("86247 RICE LG ", 1.04),
    @("727495 TEA GREEN ", 0.65),
    @("6023 COFFEE DECAF FD", 1.19),
    @("835067 PISTACHIOS SALTED", 1.99),
    @("42929 BREAD WHOLEMEAL", 0.75)
)

for ($i = 0; $i -lt $items.Length; $i++) {
    $r = $i + 2
    $name = $items[$i][0]
    $price = $items[$i][1]
    $share = -$price / 4

    $ws.Cells.Item($r, 1).Value = $name
    $ws.Cells.Item($r, 2).Value = $price

    $ws.Cells.Item($r, 3).Value = 1   # Anshuman_cont
    $ws.Cells.Item($r, 4).Value = 1   # Ritesh_cont
    $ws.Cells.Item($r, 5).Value = 1   # Rohan_cont
    $ws.Cells.Item($r, 6).Value = 1   # Ashish_cont

    $ws.Cells.Item($r, 7).Value = $share    # Anshuman_value
    $ws.Cells.Item($r, 8).Value = $share    # Ritesh_value
    $ws.Cells.Item($r, 9).Value = $share    # Rohan_value
    $ws.Cells.Item($r, 10).Value = $share   # Ashish_value
}
